# "aggiunto ricalcolo moq su ordini f non auto, sistemare il log"
#
# Work on the "nuovo" sheet (the second sheet in the workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nuovo")
$ws.Activate()

# Recalculate the MOQ: the order multiple (B8) moves from 10 to 12.
# All the downstream shared formulas in the E23:I28 block depend on
# $B$8 and will ripple through automatically.
$ws.Range("B8").Value = 12

# Add the new "moq" recalculation check below the existing table.
$ws.Range("E33").Formula = "=90/12"

# Leave the selection on E24, matching where the user ended up looking
# at the recalculated table.
$ws.Range("E24").Select() | Out-Null
